{"js": "// The document contains three Word \"fields\" (fldChar begin/instrText/fldChar end)\n// that encode M2Doc template tokens (\"m:for v | self.eClassifiers\", \"m:v.oclIsKindOf(...)\",\n// \"m:endfor\"). This edit rewrites the parser output so those fields are expressed as\n// plain literal text runs using \"{\" / \"}\" delimiters instead of Word field codes\n// (TokenIteratorFieldRewriterSplit), while keeping the surrounding formatted runs\n// (the blue \"<---Always true...\" comment) intact and reordering the oclIsKindOf(...)\n// text (with its bookmark) to sit right after \"m:v.\" inside the now-literal \"{...}\" run.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 1 = \"{m:for v | self.eClassifiers}\" field (preceded by a leading space run).\n// Paragraph 2 = \"{m:v.oclIsKindOf(ecore::EClassifier)}\" field (with the blue comment runs).\n// Paragraph 3 = \"{m:endfor}\" field.\n// Paragraph 4 = \"End of demonstration.\" (used only as the end boundary of the range).\nconst startParagraph = paragraphs.items[1];\nconst endParagraph = paragraphs.items[4];\n\nconst range = startParagraph.getRange(\"Start\").expandTo(endParagraph.getRange(\"Start\"));\n\nconst wNs = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\n\nconst newParagraphsXml =\n  \"<w:p><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r>\" +\n  \"<w:r><w:t>{m:</w:t></w:r>\" +\n  \"<w:r><w:t>for v | self.eClassifiers}</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>{</w:t></w:r>\" +\n  \"<w:r><w:t>m</w:t></w:r>\" +\n  \"<w:r><w:t>:v.</w:t></w:r>\" +\n  \"<w:r><w:t>oclIsKindOf(ecore::EClassifier</w:t></w:r>\" +\n  \"<w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/>\" +\n  \"<w:bookmarkEnd w:id=\\\"0\\\"/>\" +\n  \"<w:r><w:t>)</w:t></w:r>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\">}</w:t></w:r>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\">    </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"0000FF\\\"/><w:sz w:val=\\\"32\\\"/><w:highlight w:val=\\\"lightGray\\\"/></w:rPr><w:t>&lt;---</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:color w:val=\\\"0000FF\\\"/><w:sz w:val=\\\"32\\\"/><w:highlight w:val=\\\"lightGray\\\"/></w:rPr>\" +\n  \"<w:t>Always true:\\nNothing inferred when v (EClassifier=EClassifier) is not kind of EClassifierLiteral=EClassifier</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>{</w:t></w:r>\" +\n  \"<w:r><w:t>m:</w:t></w:r>\" +\n  \"<w:r><w:t xml:space=\\\"preserve\\\">endfor}</w:t></w:r></w:p>\";\n\nconst ooxmlPackage =\n  \"<?xml version=\\\"1.0\\\" standalone=\\\"yes\\\"?>\" +\n  \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\">\" +\n  \"<pkg:part pkg:name=\\\"/word/document.xml\\\" \" +\n  \"pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\" +\n  \"<pkg:xmlData><w:document xmlns:w=\\\"\" + wNs + \"\\\"><w:body>\" +\n  newParagraphsXml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nrange.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains three Word \"fields\" (fldChar begin/instrText/fldChar end)\n# that encode M2Doc template tokens (\"m:for v | self.eClassifiers\", \"m:v.oclIsKindOf(...)\",\n# \"m:endfor\"). This edit rewrites the parser output so those fields are expressed as\n# plain literal text runs using \"{\" / \"}\" delimiters instead of Word field codes\n# (TokenIteratorFieldRewriterSplit), while keeping the surrounding formatted runs\n# (the blue \"<---Always true...\" comment) intact and reordering the oclIsKindOf(...)\n# text (with its bookmark) to sit right after \"m:v.\" inside the now-literal \"{...}\" run.\n\n$d = $word.ActiveDocument\n\n# Paragraph 2 = \"{m:for v | self.eClassifiers}\" field (preceded by a leading space run).\n# Paragraph 3 = \"{m:v.oclIsKindOf(ecore::EClassifier)}\" field (with the blue comment runs).\n# Paragraph 4 = \"{m:endfor}\" field.\n# Paragraph 5 = \"End of demonstration.\" (used only as the end boundary of the range).\n$startParagraph = $d.Paragraphs.Item(2)\n$endParagraph = $d.Paragraphs.Item(5)\n\n$range = $d.Range($startParagraph.Range.Start, $endParagraph.Range.Start)\n\n$newParagraphsXml = (\n  '<w:p><w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>{m:</w:t></w:r>' +\n  '<w:r><w:t>for v | self.eClassifiers}</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:v.</w:t></w:r>' +\n  '<w:r><w:t>oclIsKindOf(ecore::EClassifier</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t>)</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">    </w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"0000FF\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr><w:t>&lt;---</w:t></w:r>' +\n  '<w:r><w:rPr><w:color w:val=\"0000FF\"/><w:sz w:val=\"32\"/><w:highlight w:val=\"lightGray\"/></w:rPr>' +\n  \"<w:t>Always true:`nNothing inferred when v (EClassifier=EClassifier) is not kind of EClassifierLiteral=EClassifier</w:t></w:r></w:p>\" +\n  '<w:p><w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m:</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">endfor}</w:t></w:r></w:p>'\n)\n\n$ooxmlPackage = (\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  $newParagraphsXml +\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n)\n\n$range.InsertXML($ooxmlPackage)\n"}
